$p = $ppt.ActivePresentation

# --- Slide 5: "Challenges and Successes" ---------------------------------
$s5 = $p.Slides.Item(5)

# Picture Placeholder 7 (background image) - nudge its position slightly.
# Target EMU: x=-20268, y=77653 (was x=0, y=69264). Values below are chosen
# so that PowerPoint's single-precision Left/Top round-trip to the exact
# target EMU (COM Left/Top are in points = EMU/12700).
$pic = $s5.Shapes.Item(6)
$pic.Left = -1.5959842118110237
$pic.Top  = 6.114409748818898

# TextBox 8 ("Merge Conflicts!!!" bullet list) - reposition.
# Target EMU: x=7133514, y=2832855 (was x=7315197, y=3561236).
$tb = $s5.Shapes.Item(16)
$tb.Left = 561.6940002480314
$tb.Top  = 223.05945591889764

# --- Slide 6: "Tasks and Roles" table -------------------------------------
$s6 = $p.Slides.Item(6)
$tbl = $s6.Shapes.Item(2).Table

# Header row: rename columns.
$tbl.Cell(1,2).Shape.TextFrame.TextRange.Text = "Planning/Logistics"
$tbl.Cell(1,4).Shape.TextFrame.TextRange.Text = "Styling"
$tbl.Cell(1,5).Shape.TextFrame.TextRange.Text = "Development"

# Body rows: fill in role assignments per person (columns: Planning/Logistics,
# Design, Styling, Development).
$tbl.Cell(2,2).Shape.TextFrame.TextRange.Text = "Senior"
$tbl.Cell(2,3).Shape.TextFrame.TextRange.Text = "Senior"
$tbl.Cell(2,4).Shape.TextFrame.TextRange.Text = "Lead"
$tbl.Cell(2,5).Shape.TextFrame.TextRange.Text = "Junior"

$tbl.Cell(3,2).Shape.TextFrame.TextRange.Text = "Senior"
$tbl.Cell(3,3).Shape.TextFrame.TextRange.Text = "Junior"
$tbl.Cell(3,4).Shape.TextFrame.TextRange.Text = "Junior"
$tbl.Cell(3,5).Shape.TextFrame.TextRange.Text = "Lead"

$tbl.Cell(4,2).Shape.TextFrame.TextRange.Text = "Senior"
$tbl.Cell(4,3).Shape.TextFrame.TextRange.Text = "Lead"
$tbl.Cell(4,4).Shape.TextFrame.TextRange.Text = "Senior"
$tbl.Cell(4,5).Shape.TextFrame.TextRange.Text = "Junior"

$tbl.Cell(5,2).Shape.TextFrame.TextRange.Text = "Lead"
$tbl.Cell(5,3).Shape.TextFrame.TextRange.Text = "Junior"
$tbl.Cell(5,4).Shape.TextFrame.TextRange.Text = "Junior"
$tbl.Cell(5,5).Shape.TextFrame.TextRange.Text = "Senior"
